$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 335. This shifts the existing
# rows 335-362 down to 336-363 (carrying all their original data with
# them), matching the diff where every row from 336..363 simply takes on
# the values that used to belong to the row above it.
$ws.Rows.Item(335).Insert()

# Populate the newly inserted row 335. Its contents mirror what the old
# row 335 held (same market/category/quality/etc.), except for the date
# (column D) and volume (column J), which carry new values per the diff.
$ws.Cells.Item(335, 1).Value = 10
$ws.Cells.Item(335, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(335, 3).Value = "La Araucanía"
$ws.Cells.Item(335, 4).Value = 45021
$ws.Cells.Item(335, 5).Value = 9
$ws.Cells.Item(335, 6).Value = 100112039
$ws.Cells.Item(335, 7).Value = "Ciboulette"
$ws.Cells.Item(335, 8).Value = "Sin especificar"
$ws.Cells.Item(335, 9).Value = "Primera"
$ws.Cells.Item(335, 10).Value = 45
$ws.Cells.Item(335, 11).Value = 5000
$ws.Cells.Item(335, 12).Value = 5000
$ws.Cells.Item(335, 13).Value = 5000
$ws.Cells.Item(335, 14).Value = "$/docena de atados"
$ws.Cells.Item(335, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(335, 16).Value = 1667
$ws.Cells.Item(335, 17).Value = 3
$ws.Cells.Item(335, 18).Value = "Hortaliza"
